$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column values are forced to text (NumberFormat "@") so strings such as
# "0.0650" or "35.252.20" are preserved verbatim instead of being parsed as
# numbers; the Style is then reset to "Normal" so no stray cell formatting
# is introduced (matches the original, unstyled D/E cells).

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '35.252.20'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.23%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.901.91'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  -0.16%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.727'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +9.38%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '255.26'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +3.84%  '

$ws.Range("E7").Value = '  -0.09%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '40.54'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -1.69%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.372'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +6.40%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '52.78'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.38%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0758'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +5.59%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.0987'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.62%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '2.181.51'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.21%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '12.93'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +6.85%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.724'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +4.07%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '4.95'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.52%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.922.92'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.18%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '35.231.63'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.30%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '74.60'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +3.33%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0846'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +3.38%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '243.42'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +1.16%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '13.01'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +4.74%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '5.09'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +5.50%  '

$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("E25").Value = '  +7.75%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.42'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +4.13%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '166.09'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -2.49%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '8.67'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +3.15%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '18.73'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +2.17%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.132'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +4.22%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.129.21'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +19.47%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.37'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +5.95%  '

$ws.Range("E33").Value = '  +14.59%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.64'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +23.08%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.0587'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +4.20%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '4.24'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.12%  '

$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("E38").Value = '  -2.35%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.02'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("E40").Value = '  +5.03%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '17.05'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +6.32%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '96.28'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +7.64%  '

$ws.Range("E43").Value = '  +1.83%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0650'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +4.08%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.336.36'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.06%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.43'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +2.89%  '

$ws.Range("E47").Value = '  +1.15%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '6.71'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +3.52%  '

$ws.Range("E49").Value = '  -0.55%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '45.17'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -6.23%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0753'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +6.89%  '
